$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text format for numeric-looking price cells to preserve them as text
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Set updated values
$ws.Range("D2").Value = "30.824.22"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "2.116.44"
$ws.Range("E3").Value = "  +6.73%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "333.37"
$ws.Range("E5").Value = "  +3.22%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.5321"
$ws.Range("E7").Value = "  +4.37%  "
$ws.Range("D8").Value = "0.4384"
$ws.Range("E8").Value = "  +7.00%  "
$ws.Range("D9").Value = "0.09020"
$ws.Range("E9").Value = "  +3.64%  "
$ws.Range("D10").Value = "46.13"
$ws.Range("E10").Value = "  +8.43%  "
$ws.Range("D11").Value = "1.182"
$ws.Range("E11").Value = "  +4.81%  "
$ws.Range("D12").Value = "25.06"
$ws.Range("E12").Value = "  +3.98%  "
$ws.Range("D13").Value = "2.115.32"
$ws.Range("E13").Value = "  +6.78%  "
$ws.Range("D14").Value = "6.778"
$ws.Range("E14").Value = "  +4.72%  "
$ws.Range("D15").Value = "7.839"
$ws.Range("E15").Value = "  +6.36%  "
$ws.Range("D16").Value = "97.30"
$ws.Range("E16").Value = "  +3.64%  "
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "0.00001134"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "0.06667"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").Value = "19.18"
$ws.Range("E20").Value = "  +2.21%  "
$ws.Range("D21").Value = "1.0000"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "6.350"
$ws.Range("E22").Value = "  +4.65%  "
$ws.Range("D23").Value = "30.892.73"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").Value = "12.40"
$ws.Range("E24").Value = "  +8.25%  "
$ws.Range("D25").Value = "2.362.64"
$ws.Range("E25").Value = "  +6.85%  "
$ws.Range("D26").Value = "2.271"
$ws.Range("E26").Value = "  +2.64%  "
$ws.Range("D27").Value = "22.83"
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("D28").Value = "2.588"
$ws.Range("E28").Value = "  +9.90%  "
$ws.Range("D29").Value = "163.63"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").Value = "133.73"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("D31").Value = "1.189"
$ws.Range("E31").Value = "  +4.98%  "
$ws.Range("D32").Value = "0.1081"
$ws.Range("E32").Value = "  +2.63%  "
$ws.Range("D33").Value = "6.254"
$ws.Range("E33").Value = "  +3.72%  "
$ws.Range("E34").Value = "  +5.38%  "
$ws.Range("D35").Value = "1.563"
$ws.Range("E35").Value = "  +19.16%  "
$ws.Range("D36").Value = "0.02613"
$ws.Range("E36").Value = "  +5.46%  "
$ws.Range("D37").Value = "12.92"
$ws.Range("E37").Value = "  +9.91%  "
$ws.Range("D38").Value = "5.548"
$ws.Range("E38").Value = "  +3.25%  "
$ws.Range("D39").Value = "0.06770"
$ws.Range("E39").Value = "  +4.39%  "
$ws.Range("D40").Value = "9.493"
$ws.Range("E40").Value = "  +6.45%  "
$ws.Range("D41").Value = "0.2287"
$ws.Range("E41").Value = "  +5.37%  "
$ws.Range("D42").Value = "0.6885"
$ws.Range("E42").Value = "  +5.11%  "
$ws.Range("D43").Value = "1.255"
$ws.Range("E43").Value = "  +3.14%  "
$ws.Range("D44").Value = "0.6505"
$ws.Range("E44").Value = "  +6.48%  "
$ws.Range("D45").Value = "14.12"
$ws.Range("E45").Value = "  +4.27%  "
$ws.Range("D46").Value = "0.9994"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "2.235"
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("D48").Value = "3.671"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("D49").Value = "1.277"
$ws.Range("E49").Value = "  +4.64%  "
$ws.Range("D50").Value = "83.03"
$ws.Range("E50").Value = "  +4.60%  "
$ws.Range("D51").Value = "121.82"
$ws.Range("E51").Value = "  -1.46%  "
